$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 (BePo vehicle types) ---
# New shared strings are introduced in this order: GruKw, leBefKw, FüKw, GefKw
$ws.Range("C19:K19").Value = "GruKw"
$ws.Range("L19:N19").Value = "leBefKw"
$ws.Range("O19").Value     = "FüKw"
$ws.Range("P19").Value     = "GefKw"

# --- Row 20 (counts for row 19) ---
$ws.Range("C20:K20").Value = 9
$ws.Range("L20:O20").Value = 3

# --- Row 37 (lower summary block mirrors row 20 with literal values) ---
$ws.Range("J37").Value = 9
$ws.Range("K37").Value = 3

# --- Row 12 (RD / Ortsverband vehicle types) ---
# New shared strings introduced here, after row 19's: KdoW-LNA, KdoW-OrGl
# C12:H12 stay "RTW" (unchanged)
$ws.Range("I12").Value  = "KTW"
$ws.Range("J12:K12").Value = "NEF"
$ws.Range("L12:N12").Value = "NEF"
$ws.Range("O12:P12").Value = "KdoW-LNA"
$ws.Range("Q12").Value  = "KdoW-OrGl"

# --- Row 13 (counts for row 12) ---
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 1

# --- Active cell selection ---
$ws.Range("H16").Select()
